# Apply the "3.3.1 files from RMI" update to the BUTYGV control-setting
# workbook. The only functional content change is the addition of a
# date stamp in cell C1 of the first worksheet ("About"), formatted as
# a short date (built-in numFmtId 14), giving serial value 44307
# (2021-04-21).

$wb = $excel.ActiveWorkbook

# First sheet in the workbook ("About")
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C1").Value = 44307
$ws1.Range("C1").NumberFormat = "mm-dd-yy"
